$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto data (prices + % volume changes) per the latest scrape.
# Numeric-looking price strings must be forced to Text format first so Excel
# COM does not silently coerce them to numbers (dropping trailing/leading zeros,
# or mangling tiny values into scientific notation).

$ws.Range("D2").Value = '73.005.41'
$ws.Range("E2").Value = '  +2.15%  '

$ws.Range("D3").Value = '4.001.06'
$ws.Range("E3").Value = '  +0.66%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '622.59'
$ws.Range("E5").Value = '  +15.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.46'
$ws.Range("E6").Value = '  +8.49%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.690'
$ws.Range("E7").Value = '  +0.29%  '

$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.764'
$ws.Range("E9").Value = '  +2.27%  '

$ws.Range("E10").Value = '  +0.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.72'
$ws.Range("E11").Value = '  -1.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000319'
$ws.Range("E12").Value = '  -0.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.20'
$ws.Range("E13").Value = '  +4.39%  '

$ws.Range("D14").Value = '4.632.36'
$ws.Range("E14").Value = '  +0.48%  '

$ws.Range("D15").Value = '4.002.98'
$ws.Range("E15").Value = '  +0.57%  '

$ws.Range("E16").Value = '  +7.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.19'
$ws.Range("E17").Value = '  +0.99%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.73'
$ws.Range("E18").Value = '  +0.38%  '

$ws.Range("E19").Value = '  +0.38%  '

$ws.Range("D20").Value = '72.722.33'
$ws.Range("E20").Value = '  +1.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '442.50'
$ws.Range("E21").Value = '  +2.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.12'
$ws.Range("E22").Value = '  +20.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '97.23'
$ws.Range("E23").Value = '  -0.15%  '

$ws.Range("E24").Value = '  -3.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.48'
$ws.Range("E25").Value = '  -0.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.32'
$ws.Range("E26").Value = '  +3.98%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.37'
$ws.Range("E27").Value = '  -0.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.65'
$ws.Range("E28").Value = '  -2.14%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.97'
$ws.Range("E29").Value = '  +0.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.53'
$ws.Range("E30").Value = '  -0.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.75'
$ws.Range("E31").Value = '  -2.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '14.03'
$ws.Range("E32").Value = '  +4.65%  '

$ws.Range("E33").Value = '  -0.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '72.04'
$ws.Range("E34").Value = '  +9.68%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '48.13'
$ws.Range("E35").Value = '  -6.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '641.68'
$ws.Range("E36").Value = '  -5.26%  '

$ws.Range("D37").Value = '0.0₃0900'
$ws.Range("E37").Value = '  +9.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.441'
$ws.Range("E38").Value = '  -1.10%  '

$ws.Range("E39").Value = '  -1.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.37'
$ws.Range("E40").Value = '  -1.34%  '

$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.36'
$ws.Range("E41").Value = '  +3.94%  '

$ws.Range("B42").Value = 'Dai'
$ws.Range("C42").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("E43").Value = '  +0.20%  '

$ws.Range("E44").Value = '  +1.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.71'
$ws.Range("E45").Value = '  +2.88%  '

$ws.Range("E46").Value = '  +1.18%  '

$ws.Range("E47").Value = '  -0.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.44'
$ws.Range("E48").Value = '  +2.39%  '

$ws.Range("D49").Value = '2.900.05'
$ws.Range("E49").Value = '  +10.00%  '

$ws.Range("E50").Value = '  +1.99%  '

$ws.Range("E51").Value = '  +4.10%  '
